$d = $word.ActiveDocument

function Set-CellText($cell, $text, $fontName, $fontSize) {
    $cell.Range.Text = $text
    $r = $cell.Range
    $fmtRange = $d.Range($r.Start, $r.End - 1)
    $fmtRange.Font.Name = $fontName
    $fmtRange.Font.Size = $fontSize
}

# --- Header date update ---
$d.Content.Find.Execute("2025-07-29 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-30 Wednesday", 2)

$t = $d.Tables.Item(1)

# --- Group A (row 1): simple value replacements ---
Set-CellText $t.Rows.Item(1).Cells.Item(1) "41÷7=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(1).Cells.Item(2) "49÷6=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(1).Cells.Item(3) "99÷2=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(1).Cells.Item(4) "53÷3=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(1).Cells.Item(5) "17÷7=" "TimeNewRoman" 15

# --- Group B (row 5): simple value replacements ---
Set-CellText $t.Rows.Item(5).Cells.Item(1) "64÷7=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(5).Cells.Item(2) "90÷5=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(5).Cells.Item(3) "58÷9=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(5).Cells.Item(4) "71÷7=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(5).Cells.Item(5) "15÷6=" "TimeNewRoman" 15

# --- Insert a new block (3 blank rows + 1 data row) right after Group B, ---
# --- i.e. before the row that used to be row 9 (start of Group C). ---
# --- Rows.Add(ref) inserts immediately before ref, so build back-to-front ---
# --- (data row first, then the blanks) to land in the right final order. ---
$refRow = $t.Rows.Item(9)
$newDataRow = $t.Rows.Add($refRow)
$t.Rows.Add($refRow) | Out-Null
$t.Rows.Add($refRow) | Out-Null
$t.Rows.Add($refRow) | Out-Null

Set-CellText $newDataRow.Cells.Item(1) "65÷3=" "TimeNewRoman" 15
Set-CellText $newDataRow.Cells.Item(2) "34÷9=" "TimeNewRoman" 15
Set-CellText $newDataRow.Cells.Item(3) "92÷5=" "TimeNewRoman" 15
Set-CellText $newDataRow.Cells.Item(4) "40÷4=" "TimeNewRoman" 15
Set-CellText $newDataRow.Cells.Item(5) "98÷5=" "TimeNewRoman" 15

# --- Group C (now row 13): first cell unchanged, rest replaced ---
Set-CellText $t.Rows.Item(13).Cells.Item(2) "25÷5=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(13).Cells.Item(3) "50÷8=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(13).Cells.Item(4) "22÷4=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(13).Cells.Item(5) "18÷4=" "TimeNewRoman" 15

# --- Group D (now row 17): simple value replacements ---
Set-CellText $t.Rows.Item(17).Cells.Item(1) "87÷8=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(17).Cells.Item(2) "70÷5=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(17).Cells.Item(3) "26÷6=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(17).Cells.Item(4) "24÷8=" "TimeNewRoman" 15
Set-CellText $t.Rows.Item(17).Cells.Item(5) "79÷6=" "TimeNewRoman" 15

# --- Remove the trailing block (old Group E: data row + 3 blanks), now rows 21-24 ---
$t.Rows.Item(21).Delete()
$t.Rows.Item(21).Delete()
$t.Rows.Item(21).Delete()
$t.Rows.Item(21).Delete()
